$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").Value = ""

$ws.Range("H29").Value = 8999
$ws.Range("J29").Value = 8999
$ws.Range("L29").Value = 26997
$ws.Range("N29").Value = -27559

$ws.Range("H41").Value = 868.875
$ws.Range("I41").Value = 664.2857
$ws.Range("J41").Value = 1155.3
$ws.Range("K41").Value = 664.2857
$ws.Range("L41").Value = 1155.3
$ws.Range("M41").Value = -224.2857
$ws.Range("N41").Value = -2035.3

$ws.Range("H51").Value = 7342.1
$ws.Range("I51").Value = 5854.75
$ws.Range("K51").Value = 5854.75
$ws.Range("M51").Value = -5370.75

$ws.Range("H80").Value = 4160.3335
$ws.Range("I80").Value = 2455.2
$ws.Range("K80").Value = 7365.599999999999
$ws.Range("M80").Value = -6367.599999999999

$ws.Range("H83").Value = 4160.3335
$ws.Range("I83").Value = 2455.2
$ws.Range("K83").Value = 22096.8
$ws.Range("M83").Value = -17104.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9423.725
$ws.Range("I32").Value = 6856.523
$ws.Range("J32").Value = 17492.072
$ws.Range("K32").Value = 6856.523
$ws.Range("L32").Value = 17492.072
$ws.Range("M32").Value = -6569.523
$ws.Range("N32").Value = -18066.072

$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("N68").Value = 0

$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("N71").Value = 0

$ws.Range("H74").Value = 1261.4546
$ws.Range("I74").Value = 887.6
$ws.Range("K74").Value = 887.6
$ws.Range("M74").Value = -13.60000000000002

$ws.Range("H77").Value = 1261.4546
$ws.Range("I77").Value = 887.6
$ws.Range("K77").Value = 4438
$ws.Range("M77").Value = -70

$ws.Range("H97").Value = 1873.9524
$ws.Range("I97").Value = 1860.3334
$ws.Range("J97").Value = 1908
$ws.Range("K97").Value = 1860.3334
$ws.Range("L97").Value = 1908
$ws.Range("M97").Value = -1364.3334
$ws.Range("N97").Value = -2900

$ws.Range("H102").Value = 2128.611
$ws.Range("I102").Value = 1900.2354
$ws.Range("K102").Value = 1900.2354
$ws.Range("M102").Value = -278.2354

$ws.Range("H132").Value = 3502.3076
$ws.Range("I132").Value = 2764.4707
$ws.Range("J132").Value = 4896
$ws.Range("K132").Value = 8293.4121
$ws.Range("L132").Value = 14688
$ws.Range("M132").Value = -5763.4121
$ws.Range("N132").Value = -19748

$ws.Range("H135").Value = 51664.332
$ws.Range("J135").Value = 51664.332
$ws.Range("L135").Value = 51664.332
$ws.Range("N135").Value = -61804.332

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2882.5833
$ws.Range("I86").Value = 1371.5
$ws.Range("J86").Value = 7415.8335
$ws.Range("K86").Value = 1371.5
$ws.Range("L86").Value = 7415.8335
$ws.Range("M86").Value = -248.5
$ws.Range("N86").Value = -9661.833500000001

$ws.Range("H89").Value = 2882.5833
$ws.Range("I89").Value = 1371.5
$ws.Range("J89").Value = 7415.8335
$ws.Range("K89").Value = 6857.5
$ws.Range("L89").Value = 37079.1675
$ws.Range("M89").Value = -1241.5
$ws.Range("N89").Value = -48311.1675

$ws.Range("H99").Value = 13870.743
$ws.Range("I99").Value = 15635.2
$ws.Range("J99").Value = 3284
$ws.Range("K99").Value = 15635.2
$ws.Range("L99").Value = 3284
$ws.Range("M99").Value = -14137.2
$ws.Range("N99").Value = -6280

$ws.Range("H105").Value = 3287.7144
$ws.Range("I105").Value = 3237.1
$ws.Range("K105").Value = 3237.1
$ws.Range("M105").Value = -1490.1

$ws.Range("H107").Value = 1785
$ws.Range("I107").Value = 1749.1666
$ws.Range("K107").Value = 1749.1666
$ws.Range("M107").Value = 170.8334

$ws.Range("H134").Value = 4523.773
$ws.Range("I134").Value = 3379.2144
$ws.Range("J134").Value = 6526.75
$ws.Range("K134").Value = 10137.6432
$ws.Range("L134").Value = 19580.25
$ws.Range("M134").Value = -7602.643199999999
$ws.Range("N134").Value = -24650.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1815.5312
$ws.Range("I58").Value = 1073.1904
$ws.Range("J58").Value = 3232.7273
$ws.Range("K58").Value = 1073.1904
$ws.Range("L58").Value = 3232.7273
$ws.Range("M58").Value = -870.1904
$ws.Range("N58").Value = -3638.7273

$ws.Range("H62").Value = 5986.8887
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").Value = ""

$ws.Range("H65").Value = 5986.8887
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").Value = ""

$ws.Range("H109").Value = 33252.31
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 33252.31
$ws.Range("K109").Value = 0
$ws.Range("M109").Value = 33252.31
$ws.Range("N109").Value = -35332.31

$ws.Range("H132").Value = 2580.48
$ws.Range("I132").Value = 2580.48
$ws.Range("K132").Value = 7741.440000000001
$ws.Range("M132").Value = -5211.440000000001

$ws.Range("H136").Value = 1815.5312
$ws.Range("I136").Value = 1073.1904
$ws.Range("J136").Value = 3232.7273
$ws.Range("K136").Value = 3219.5712
$ws.Range("L136").Value = 9698.1819
$ws.Range("M136").Value = -669.5711999999999
$ws.Range("N136").Value = -14798.1819

$ws.Range("H141").Value = 92821.42999999999
$ws.Range("J141").Value = 92821.42999999999
$ws.Range("L141").Value = 92821.42999999999
$ws.Range("N141").Value = -103181.43

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 54726.5
$ws.Range("I5").Value = 10004
$ws.Range("J5").Value = 99449
$ws.Range("K5").Value = 30012
$ws.Range("L5").Value = 298347
$ws.Range("M5").Value = -29900
$ws.Range("N5").Value = -298571

$ws.Range("H11").Value = 1238.6666
$ws.Range("I11").Value = 610
$ws.Range("J11").Value = 1867.3334
$ws.Range("K11").Value = 1830
$ws.Range("L11").Value = 5602.0002
$ws.Range("M11").Value = -1690
$ws.Range("N11").Value = -5882.0002

$ws.Range("H131").Value = 8858
$ws.Range("I131").Value = 5241.75
$ws.Range("J131").Value = 10173
$ws.Range("K131").Value = 15725.25
$ws.Range("L131").Value = 30519
$ws.Range("M131").Value = -10685.25
$ws.Range("N131").Value = -40599

$ws.Range("H132").Value = 1200
$ws.Range("J132").Value = 1200
$ws.Range("L132").Value = 10800
$ws.Range("N132").Value = -15860

$ws.Range("H135").Value = 54726.5
$ws.Range("I135").Value = 10004
$ws.Range("J135").Value = 99449
$ws.Range("K135").Value = 90036
$ws.Range("L135").Value = 895041
$ws.Range("M135").Value = -87501
$ws.Range("N135").Value = -900111

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 4166931.5
$ws.Range("I2").Value = 243.76923
$ws.Range("K2").Value = 243.76923
$ws.Range("M2").Value = -130.76923

$ws.Range("H70").Value = 127669.11
$ws.Range("I70").Value = 225943.2
$ws.Range("J70").Value = 4826.5
$ws.Range("K70").Value = 225943.2
$ws.Range("L70").Value = 4826.5
$ws.Range("M70").Value = -225673.2
$ws.Range("N70").Value = -5366.5

$ws.Range("H73").Value = 127669.11
$ws.Range("I73").Value = 225943.2
$ws.Range("J73").Value = 4826.5
$ws.Range("K73").Value = 225943.2
$ws.Range("L73").Value = 4826.5
$ws.Range("M73").Value = -225007.2
$ws.Range("N73").Value = -6698.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7300
$ws.Range("I7").Value = 8872.444
$ws.Range("K7").Value = 8872.444
$ws.Range("M7").Value = -8760.444

$ws.Range("H16").Value = 1662.4117
$ws.Range("I16").Value = 1876.7142
$ws.Range("K16").Value = 1876.7142
$ws.Range("M16").Value = -1706.7142

$ws.Range("H40").Value = 15937.833
$ws.Range("I40").Value = 52502
$ws.Range("J40").Value = 8625
$ws.Range("K40").Value = 52502
$ws.Range("L40").Value = 8625
$ws.Range("M40").Value = -52366
$ws.Range("N40").Value = -8897

$ws.Range("H126").Value = 7300
$ws.Range("I126").Value = 8872.444
$ws.Range("K126").Value = 26617.332
$ws.Range("M126").Value = -24147.332

$ws.Range("H135").Value = 80000
$ws.Range("J135").Value = 80000
$ws.Range("L135").Value = 80000
$ws.Range("N135").Value = -90140

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("M15").Value = ""

$ws.Range("H22").Value = 2500
$ws.Range("J22").Value = 2500
$ws.Range("L22").Value = 2500
$ws.Range("N22").Value = -3086

$ws.Range("H69").Value = 22499.834

$ws.Range("H72").Value = 22499.834

$ws.Range("H107").Value = 749.6667
$ws.Range("J107").Value = 500
$ws.Range("L107").Value = 1500
$ws.Range("N107").Value = -5340

$ws.Range("H109").Value = 31562.5
$ws.Range("J109").Value = 31562.5
$ws.Range("L109").Value = 31562.5
$ws.Range("N109").Value = -34336.5

$ws.Range("H126").Value = 1638.2
$ws.Range("I126").Value = 1561.909
$ws.Range("J126").Value = 2197.6667
$ws.Range("K126").Value = 4685.727000000001
$ws.Range("L126").Value = 6593.000100000001
$ws.Range("M126").Value = -2215.727000000001
$ws.Range("N126").Value = -11533.0001
